$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The pipeline's landing/raw/control databases for "ad_works_dw" are being
# reconfigured with a "yetl_" naming prefix. Rename every occurrence
# (database name cells, and the fully-qualified "landing.<db>.<table>" /
# "raw.<db>.*" depends_on references that are built from those names).
$ws.Cells.Replace("landing_ad_works_dw", "yetl_landing_ad_works_dw")
$ws.Cells.Replace("raw_ad_works_dw", "yetl_raw_ad_works_dw")
$ws.Cells.Replace("control_ad_works_dw", "yetl_control_ad_works_dw")

# Move the active selection to C4 (matches the saved cursor position).
$ws.Range("C4").Select()
